# Apply cryptos.xlsx price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.617.86'
$ws.Range("E2").Value = '  -2.29%  '

$ws.Range("D3").Value = '2.896.79'
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.09'
$ws.Range("E5").Value = '  -4.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.59'
$ws.Range("E6").Value = '  -3.35%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.44%  '

$ws.Range("D9").Value = '2.895.21'
$ws.Range("E9").Value = '  -1.90%  '

$ws.Range("E10").Value = '  -1.68%  '

$ws.Range("E11").Value = '  -2.16%  '

$ws.Range("E12").Value = '  -2.20%  '

$ws.Range("E13").Value = '  -0.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.85'
$ws.Range("E14").Value = '  -2.87%  '

$ws.Range("E15").Value = '  -0.61%  '

$ws.Range("D16").Value = '3.376.48'
$ws.Range("E16").Value = '  -1.91%  '

$ws.Range("D17").Value = '61.604.87'
$ws.Range("E17").Value = '  -2.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.55'
$ws.Range("E18").Value = '  -2.16%  '

$ws.Range("D19").Value = '2.902.91'
$ws.Range("E19").Value = '  -1.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '433.60'
$ws.Range("E20").Value = '  -1.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.08'
$ws.Range("E21").Value = '  -3.19%  '

$ws.Range("E22").Value = '  -1.58%  '

$ws.Range("E23").Value = '  -2.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.37'
$ws.Range("E24").Value = '  -2.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.94'
$ws.Range("E25").Value = '  +1.42%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.95'
$ws.Range("E27").Value = '  -10.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.02'
$ws.Range("E28").Value = '  -5.44%  '

$ws.Range("E29").Value = '  +5.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.01'
$ws.Range("E30").Value = '  -4.18%  '

$ws.Range("E31").Value = '  -4.04%  '

$ws.Range("E32").Value = '  -7.22%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("E34").Value = '  -1.68%  '

$ws.Range("E35").Value = '  -3.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.960'
$ws.Range("E36").Value = '  -3.19%  '

$ws.Range("E37").Value = '  -3.72%  '

$ws.Range("E38").Value = '  -1.68%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  -5.28%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.82'
$ws.Range("E40").Value = '  -8.59%  '

$ws.Range("E41").Value = '  -3.43%  '

$ws.Range("E42").Value = '  -2.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.64'
$ws.Range("E43").Value = '  +2.29%  '

$ws.Range("E44").Value = '  -5.01%  '

$ws.Range("D45").Value = '2.692.94'
$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.86'
$ws.Range("E46").Value = '  -1.99%  '

$ws.Range("E47").Value = '  -0.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '348.66'
$ws.Range("E48").Value = '  -3.15%  '

$ws.Range("E50").Value = '  -1.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.61'
$ws.Range("E51").Value = '  -5.04%  '
